# OverallStatus.xlsx - "Updated current status of UDQ screens"
#
# The UDQ status grid on Sheet1 (columns D/E/F, rows 5-38) holds one of three
# states per screen: "Done" (green text), "Not Done" (red/orange text) and
# "In Progress" (red/orange text on yellow fill). Several screens moved from
# "Done" to "In Progress" or "Not Done". The I/J helper formulas (COUNTIF)
# and the two charts that are bound to them recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Template cells that already carry the exact target formatting, used as a
# format source for PasteSpecial so the copied cell's font / fill match the
# existing "Done" / "Not Done" / "In Progress" look exactly.
$doneTemplate       = $ws.Range("D5")
$notDoneTemplate    = $ws.Range("F6")
$inProgressTemplate = $ws.Range("D6")

function Set-Status {
    param(
        [string]$CellAddress,
        [string]$Status
    )

    switch ($Status) {
        "Done"        { $template = $doneTemplate }
        "Not Done"    { $template = $notDoneTemplate }
        "In Progress" { $template = $inProgressTemplate }
    }

    $target = $ws.Range($CellAddress)
    $template.Copy()
    $target.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $target.Value2 = $Status
}

$excel.CutCopyMode = 0

# Row 7 - Taxability for Authority
Set-Status "D7" "In Progress"
Set-Status "E7" "In Progress"
Set-Status "F7" "Not Done"

# Row 8 - Companies
Set-Status "E8" "In Progress"

# Row 9 - Garnishment Disposable Wages
Set-Status "D9" "In Progress"
Set-Status "E9" "In Progress"
Set-Status "F9" "Not Done"

# Row 16 - Custom Taxability for Authority
Set-Status "F16" "Not Done"

# Row 21 - Garnishment Parameters
Set-Status "F21" "Not Done"

# Row 25 - Garnishments Requiring Filing Status
Set-Status "D25" "In Progress"
Set-Status "E25" "In Progress"
Set-Status "F25" "Not Done"

# Rows 33-37 - All Mapped Pay Codes / Tax Codes / Tax Types / Experience Rates / Populated V3 States
Set-Status "D33" "In Progress"
Set-Status "E33" "In Progress"
Set-Status "D34" "In Progress"
Set-Status "E34" "In Progress"
Set-Status "D35" "In Progress"
Set-Status "E35" "In Progress"
Set-Status "D36" "In Progress"
Set-Status "E36" "In Progress"
Set-Status "D37" "In Progress"
Set-Status "E37" "In Progress"

$excel.CutCopyMode = 0

# Reflect the last-edited cell as the active selection, matching the saved
# workbook view (scrolled back to the top, E8 selected).
$ws.Activate()
$ws.Range("E8").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1

$wb.Save()
